$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove existing hyperlinks on column F before rewriting values/links
$ws.Hyperlinks.Delete()

# Row 2
$ws.Range("A2").Value = "11/08/2023 11:01"
$ws.Range("B2").Value = "Las lentejas también son para el verano: esta ensalada es fácil, saludable, rápida y muy rica"
$ws.Range("C2").Value = "Inés Vazquez Noya"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = "https://www.directoalpaladar.com/recetas-vegetarianas/lentejas-tambien-para-verano-esta-ensalada-facil-saludable-rapida-muy-rica"

# Row 3
$ws.Range("A3").Value = "11/08/2023 10:01"
$ws.Range("B3").Value = "Espaguetis con salsa de berenjena, tomate y queso feta: una receta casera tan deliciosa que la querrás hacer una y otra vez"
$ws.Range("C3").Value = "Carmen Tía Alia"
$ws.Range("D3").Value = "Espaguetis (o cualquier otra pasta)300 gCebolla 50 gDiente de ajo 2 Berenjena 200 gVino blanco 50 mlTomates en conserva natural 350 gTomate concentrado 30 gHierbas provenzales 1 cucharaditaOrégano seco 1 cucharaditaQueso feta desmigado75 gSal Pimienta negra molida Aceite de oliva virgen extra Azúcar (opcional)"
$ws.Range("E3").Value = "Picamos finamente la cebolla y rallamos los dientes de ajo. Calentamos un fondo de aceite de oliva virgen extra en una sartén amplia y los pochamos durante diez minutos. Mientras tanto, pelamos y cortamos la berenjena en pequeños dados. La incorporamos a la sartén, sazonamos al gusto y continuamos sofriendo las verduras. Esta vez necesitaremos más tiempo, unos 15-20 minutos, para que la berenjena quede bien tierna y melosa. Añadimos el vino y cocemos a fuego fuerte hasta que se haya evaporado.A continuación agregamos los tomates, troceamos o machacados, junto con el tomate concentrado y las hierbas. Ajustamos el punto de sal y, si lo deseamos, agregamos una pizca de azúcar para equilibrar la acidez. Cocemos a fuego muy lento durante 45-50 minutos o hasta obtener una salsa espesa.Cocemos la pasta en abundante agua con sal y siguiendo las indicaciones del fabricante. Podemos usar parte de este agua para aligerar un poco la salsa, en caso de ser necesario. Agregamos el queso feta desmigado a la salsa, removemos y listo. Ya podemos servir nuestros espaguetis con la salsa por encima y previamente mezclado."
$ws.Range("F3").Value = "https://www.directoalpaladar.com/recetas-vegetarianas/espaguetis-salsa-berenjena-tomate-queso-feta-receta-casera-deliciosa-que-querras-hacer-otra-vez"

# Row 4
$ws.Range("A4").Value = "10/08/2023 16:01"
$ws.Range("B4").Value = "Hamburguesa de garbanzos: cómo resolver una cena con esta receta fácil y rápida para todos los públicos"
$ws.Range("C4").Value = "Inés Vazquez Noya"
$ws.Range("D4").Value = "Garbanzos cocidos en conserva 200 gCebolla blanca 0.5 Semillas de sésamo cucharadas2 Comino molido 10 gPimentón ahumado 10 gOrégano seco o fresco15 gSal y pimienta al gusto Aceite de oliva virgen extra 2 cucharadas para la hamburguesa y un extra para cocinar en la sartén"
$ws.Range("E4").Value = "En un robot de cocina, batidora o procesador de alimentos, volcar los garbanzos cocidos (y limpios de la conserva) junto a la cebolla, las semillas de sésamo, los condimentos y el aceite de oliva. Procesar unos minutos nomás, de 2 a 3 minutos, para que se integren sin alcanzar una consistencia de puré homogéneo. Si resulta muy húmeda la mezcla, incorporar unas cucharadas de harina, ya que se debe poder manejar con las manos.  Con las manos, levemente húmedas, formar las hamburguesas dándole la forma deseada. Saldrán dos unidades aproximadamente. En una sartén con aceite de oliva o antiadherente a fuego medio-bajo, cocinar los medallones hasta que estén doradas de cada lado, volteando con mucho cuidado para que no se rompan. Unos 3 minutos de cada lado. Reservar mientras se preparan los ingredientes restantes para montar la hamburguesa; tomate en rodajas y lechugas varias. Pero puedes complementar con tus ingredientes preferidos."
$ws.Range("F4").Value = "https://www.directoalpaladar.com/recetas-vegetarianas/hamburguesa-garbanzos-bote-legumbres-tienes-comida-bien-rica"

# Row 5
$ws.Range("A5").Value = "10/08/2023 10:01"
$ws.Range("B5").Value = "Arroz al coco, un truco para darle más sabor al cereal de todos los días"
$ws.Range("C5").Value = "Inés Vazquez Noya"
$ws.Range("D5").Value = "Arroz de grano largo o Jazmín o basmati150 gLeche de coco 250 mlAgua 250 mlOrégano seco o fresco - opcionalCoco rallado 3 a 4 cucharadas"
$ws.Range("E5").Value = "En una olla, volcar el arroz. Agregar la leche de coco y el agua. Revolver con cuchara de madera para diluir e integrar los ingredientes. Llevar a ebullición, a fuego medio, y dejar cocinar por 20 minutos aproximadamente. Unos minutos antes de que el liquido se haya consumido completamente, espolvorear el coco rallado sobre el cereal. Continuar la cocción unos 5 minutos más. Una vez listo, perfumar con perejil fresco o seco si gustan."
$ws.Range("F5").Value = "https://www.directoalpaladar.com/recetas-vegetarianas/arroz-al-coco-truco-para-darle-sabor-al-cereal-todos-dias"

# Row 6
$ws.Range("A6").Value = "09/08/2023 16:00"
$ws.Range("B6").Value = "Palomitas de maíz al curry, un picoteo diferente y adictivo"
$ws.Range("C6").Value = "Inés Vazquez Noya"
$ws.Range("D6").Value = "Maíz pisingallo especial para palomitas60 gRalladura de limón 1 Curry molido 12 gSirope de ágave o miel30 mlSal en escamas pizca generosa"
$ws.Range("E6").Value = "Primero, será preparar las palomitas. En un cazo (o una sartén) sin aceite poner a calentar a fuego medio. Volcar los granos de maíz y tapar hasta que empiecen a explotar. Remover con cuidado para que todos los granos tengan contacto con el fondo de la olla y el calor para explotar. Esperar hasta que se hagan todas las palomitas (aproximadamente unos 4 minutos).Una vez que dejen de explotar, retirar del fuego y salar inmediatamente con la sal en escamas para que se adhiera con el calor residual.  Volcar luego sobre una placa para el horno con papel vegetal o silpat, y acomodar por toda la superficie. Espolvorear el curry molido, la ralladura de limón y el sirope de ágave. Remover constantemente con una espátula para adherir los ingredientes a las palomitas. Hornear por 10 minutos a 180 ºC para que tome un color dorado. Retirar del horno con cuidado y dejar reposar 10 minutos antes de picar."
$ws.Range("F6").Value = "https://www.directoalpaladar.com/recetas-vegetarianas/palomitas-maiz-al-curry-picoteo-diferente-adictivo"

# Row 7
$ws.Range("A7").Value = "07/08/2023 16:00"
$ws.Range("B7").Value = "Tallarines de pepino con yogur, una receta que simplemente grita verano"
$ws.Range("C7").Value = "Inés Vazquez Noya"
$ws.Range("D7").Value = "Pepino 2 Yogur natural bote pequeño1 Ralladura de limón 0.5 Sal y pimienta al gusto Semillas de sésamo Eneldo seco o hierba de preferencia"
$ws.Range("E7").Value = "Para los tallarines de pepino, lavar y secar bien las unidades. Hacer los espirales usando un espiralizador, una mandolina, un cortador manual o un pelador con dientes. Una vez listos, reservar en un cuenco hondo. Agregar la ralladura de limón y el yogur natural. Revolver con cuidado para no romper los tallarines. Condimentar con la sal y la pimienta, y las semillas de sésamo. Conservar en la nevera hasta el momento de servir y espolvorear el eneldo seco, o cualquier hierba fresca de preferencia."
$ws.Range("F7").Value = "https://www.directoalpaladar.com/recetas-vegetarianas/tallarines-pepino-yogur-receta-que-simplemente-grita-verano"

# Row 8
$ws.Range("A8").Value = "07/08/2023 09:00"
$ws.Range("B8").Value = "Recetas frescas para disfrutar el verano en el menú semanal vegetariano del 7 de agosto"
$ws.Range("C8").Value = "Inés Vazquez Noya"
$ws.Range("D8").Value = ""
$ws.Range("E8").Value = ""
$ws.Range("F8").Value = "https://www.directoalpaladar.com/recetas-vegetarianas/recetas-frescas-para-disfrutar-verano-menu-semanal-vegetariano-7-agosto"

# Row 9
$ws.Range("A9").Value = "07/08/2023 06:00"
$ws.Range("B9").Value = "La `"paradoja del queso`": por qué los vegetarianos no quieren comer carne, pero les dan igual los lácteos"
$ws.Range("C9").Value = "Liliana Fuchs"
$ws.Range("D9").Value = ""
$ws.Range("E9").Value = ""
$ws.Range("F9").Value = "https://www.directoalpaladar.com/recetas-vegetarianas/paradoja-queso-que-vegetarianos-no-quieren-comer-carne-les-dan-igual-lacteos"

# Row 10
$ws.Range("A10").Value = "06/08/2023 17:00"
$ws.Range("B10").Value = "Paté rosa de remolacha y tofu, una receta para decorar la mesa"
$ws.Range("C10").Value = "Inés Vazquez Noya"
$ws.Range("D10").Value = "Tofu blando 100 gRemolacha cocida 50 gAgua de la remolacha - 2 cucharadasComino molido 10 gSal y pimienta al gusto"
$ws.Range("E10").Value = "Primero, tomar el tofu y quitar el líquido del envase. Presionar encima con un plato, de tal forma que vaya soltando todo el líquido y luego escurrir entre dos hojas de papel de cocina. Cortar en cubos o romper con las manos. En un robot de cocina, agregar las porciones de tofu, junto a la remolacha troceada, el comino molido, la sal, la pimienta y la cucharada de agua de la remolacha. Triturar hasta obtener una pasta homogénea.  De ser necesario, agregar 1 a 2 cucharaditas más de agua, si ves que la preparación necesita mas hidratación hasta alcanzar una textura cremosa."
$ws.Range("F10").Value = "https://www.directoalpaladar.com/recetas-vegetarianas/pate-rosa-remolacha-tofu-receta-para-decorar-mesa"

# Row 11
$ws.Range("A11").Value = "04/08/2023 15:01"
$ws.Range("B11").Value = "Este desayuno con huevos, típico de Oriente Medio, es la cena veraniega perfecta"
$ws.Range("C11").Value = "Miguel Ayuso Rejas"
$ws.Range("D11").Value = ""
$ws.Range("E11").Value = ""
$ws.Range("F11").Value = "https://www.directoalpaladar.com/recetas-vegetarianas/este-desayuno-huevos-tipico-oriente-medio-cena-veraniega-perfecta"

# Row 12
$ws.Range("A12").Value = "03/08/2023 17:01"
$ws.Range("B12").Value = "Esta comida callejera italiana, sabrosa y fácil, es mi receta favorita con calabacín"
$ws.Range("C12").Value = "Miguel Ayuso Rejas"
$ws.Range("D12").Value = ""
$ws.Range("E12").Value = ""
$ws.Range("F12").Value = "https://www.directoalpaladar.com/recetas-vegetarianas/esta-comida-callejera-italiana-sabrosa-facil-mi-receta-favorita-calabacin"

# Row 13
$ws.Range("A13").Value = "03/08/2023 16:01"
$ws.Range("B13").Value = "Muffins salados de boniato, una receta para chuparse los dedos"
$ws.Range("C13").Value = "Inés Vazquez Noya"
$ws.Range("D13").Value = "Harina integral 200 gLevadura química cucharadita1 Sal y pimienta al gusto Boniato cocido161 gLeche o bebida vegetal de almendras, de avena o de soja300 mlAceite de girasol 30 mlEspinaca fresca 50 gGuisantes 30 gComino molido 1 cuchara aproximadamente10 g"
$ws.Range("E13").Value = "Encender el horno a 190 ºC.En un cuenco mezclar los ingredientes secos: harina integral, levadura química, sal y pimienta. Integrar con un batidor de varilla o un tenedor. Reservar. Por otra parte, en un robot de cocina o licuadora, triturar 100 g de boniato ya cocido con la leche o bebida y aceite de girasol. Pulsar por varios segundos hasta conseguir una bebida liquida de color naranja. Agregar este último al recipiente con la harina e integrar hasta tener una masa homogénea. Momento de agregar la espinaca fresca, los guisantes y los 80 g de boniato restantes en cubitos, más el comino molido. Revolver nuevamente.  Acomodar en molde de muffins de silicona y verter la mezcla anterior sin alcanzar la superficie. De utilizar moldes de acero, pintar con aceite toda la superficie. Cocinar por 30 minutos. Dejar reposar unos 10 minutos antes de desmoldar y disfrutar."
$ws.Range("F13").Value = "https://www.directoalpaladar.com/recetas-vegetarianas/muffins-salados-boniato-receta-para-chuparse-dedos"

# Row 14
$ws.Range("A14").Value = "02/08/2023 11:01"
$ws.Range("B14").Value = "Refrescante, saciante y sin encender un fuego: esta receta de calabacín es una cena saludable de verano perfecta"
$ws.Range("C14").Value = "Liliana Fuchs"
$ws.Range("D14").Value = ""
$ws.Range("E14").Value = ""
$ws.Range("F14").Value = "https://www.directoalpaladar.com/recetas-vegetarianas/refrescante-saciante-encender-fuego-esta-receta-calabacin-cena-saludable-verano-perfecta"

# Row 15
$ws.Range("A15").Value = "02/08/2023 07:04"
$ws.Range("B15").Value = "Así son los nuevos postres veganos de Dia (que no son yogures) con tres bases vegetales diferentes"
$ws.Range("C15").Value = "Liliana Fuchs"
$ws.Range("D15").Value = ""
$ws.Range("E15").Value = ""
$ws.Range("F15").Value = "https://www.directoalpaladar.com/recetas-vegetarianas/asi-nuevos-postres-veganos-dia-que-no-yogures-tres-bases-vegetales-diferentes"

# Row 16
$ws.Range("A16").Value = "31/07/2023 09:00"
$ws.Range("B16").Value = "Despedimos el séptimo mes del año en el menú semanal vegetariano del 31 de julio"
$ws.Range("C16").Value = "Inés Vazquez Noya"
$ws.Range("D16").Value = ""
$ws.Range("E16").Value = ""
$ws.Range("F16").Value = "https://www.directoalpaladar.com/recetas-vegetarianas/despedimos-septimo-mes-ano-menu-semanal-vegetariano-31-julio"

# Row 17
$ws.Range("A17").Value = "30/07/2023 16:00"
$ws.Range("B17").Value = "Palitos de queso, la receta más fácil y rápida del aperitivo del verano"
$ws.Range("C17").Value = "Carmen Tía Alia"
$ws.Range("D17").Value = "Queso halloumi 200 gMiel 15 mlSemillas de sésamo 5 g"
$ws.Range("E17").Value = "Cortamos el queso halloumi en diez piezas de un tamaño aproximado y con forma de lingote. Los ponemos en remojo con abundante agua y los dejamos reposar en la nevera 30-40 minutos. Con esto conseguiremos rebajar un poco el punto de sal que tiene este queso. Secamos los palitos de halloumi con papel absorbente y los ensartamos en brochetas, aunque también los podemos dejar tal cual. Los marcamos en una parrilla bien caliente por todas sus caras. Después los pincelamos con miel y los espolvoreamos con semillas de sésamo. Servimos inmediatamente."
$ws.Range("F17").Value = "https://www.directoalpaladar.com/recetas-vegetarianas/palitos-queso-receta-facil-rapida-aperitivo-verano"

# Row 18
$ws.Range("A18").Value = "25/07/2023 16:00"
$ws.Range("B18").Value = "Galletas de naranja con harina de garbanzos y avena, una receta fácil para acompañar la merienda"
$ws.Range("C18").Value = "Inés Vazquez Noya"
$ws.Range("D18").Value = "Harina de garbanzo 60 gAvena molida30 gHarina integral 30 gCoco rallado 20 gMantequilla de cacahuete natural 20 mlNaranja zumo y ralladura1 Azúcar moreno 10 gEsencia de vainilla cucharadita1 Agua cantidad necesariaAceite de girasol 20 mlSal pizcaCanela molida 1 cucharadita (opcional)"
$ws.Range("E18").Value = "Encender el horno a 180 ºC.En un cuenco, mezclar los ingredientes secos: la harina de garbanzos, la avena molida, la harina integral, el coco rallado, el azúcar, la levadura química y la ralladura de naranja. Integrar. Agregar luego los líquidos: el aceite de girasol, el zumo de naranja, la esencia de vainilla, la mantequilla de cacahuete y el agua. Sumar una pizca de sal para que realce el sabor. Con ayuda de una cuchara de madera mezclar y terminar de integrar todo con la mano.Con las manos o una cucharada, formar 8 a 10 bolas con la masa y acomodar en una bandeja para el horno con papel de cocina o silpat. De acuerdo al tamaño, será la cantidad de galletas. Aplanar cada bola con un tenedor untado en harina o levemente húmedo. Cocinar por 10 minutos hasta que estén ligeramente doradas en la base. Dejar las galletas en la bandeja durante un minuto y transferir a una rejilla para enfriar antes de servir, así tendrán más cuerpo."
$ws.Range("F18").Value = "https://www.directoalpaladar.com/recetas-vegetarianas/galletas-naranja-harina-garbanzos-avena-receta-facil-para-acompanar-merienda"

# Row 19
$ws.Range("A19").Value = "24/07/2023 10:12"
$ws.Range("B19").Value = "Receta de ensalada coreana de patata o gamja, tan cremosa y refrescante que resulta adictiva"
$ws.Range("C19").Value = "Carmen Tía Alia"
$ws.Range("D19").Value = "Patata mediana2 Zanahoria 1 Pepino 1 Cebolla morada 1 Huevo 2 Maíz en grano (opcional)50 gMayonesa 3 cucharadasVinagre de arroz 2 cucharadasAzúcar 1 cucharadaMostaza de Dijon 1 cucharadita"
$ws.Range("E19").Value = "Pelamos las patatas y las cortamos en trozos de igual tamaño. Calentamos abundante agua con sal en una cacerola y, cuando arranque a hervir, los introducimos junto con los huevos. Retiramos estos últimos después de diez minutos y continuamos cociendo la patata hasta que esté tierna. Pelamos el pepino, la zanahoria y la cebolla morada y picamos todo finamente. Introducimos en un bol y añadimos un poco de sal, removemos y dejamos que se reblandezcan mientras sueltan sus jugos. Pelamos también los huevos y los picamos muy fino.Preparamos la salsa mezclando la mayonesa con el vinagre de arroz, el azúcar y la mostaza. Cuando la patata esté tierna la escurrimos y la pasamos por un pasapurés. Agregamos los huevos duros, las verduras -bien escurridas- y el maíz (opcional).Incorporamos la salsa y removemos para integrar todo bien. Probamos el punto de sal y añadimos más si lo consideramos necesario. Dejamos reposar la ensalada en la nevera, donde la patata irá absorbiendo la salsa poco a poco. Cuando esté bien fría, ya está lista para servir."
$ws.Range("F19").Value = "https://www.directoalpaladar.com/recetas-vegetarianas/receta-ensalada-coreana-patata-gamja-cremosa-refrescante-que-resulta-adictiva"

# Row 20
$ws.Range("A20").Value = "24/07/2023 09:00"
$ws.Range("B20").Value = "Recetas ligerísimas para disfrutar el verano en el menú semanal vegetariano del 24 de julio"
$ws.Range("C20").Value = "Inés Vazquez Noya"
$ws.Range("D20").Value = ""
$ws.Range("E20").Value = ""
$ws.Range("F20").Value = "https://www.directoalpaladar.com/recetas-vegetarianas/recetas-ligerisimas-para-disfrutar-verano-menu-semanal-vegetariano-24-julio"

# Row 21
$ws.Range("A21").Value = "24/07/2023 07:00"
$ws.Range("B21").Value = "No sabía que mis ensaladas de legumbres tenían un fallo. Todo ha cambiado desde que no me salto este truco"
$ws.Range("C21").Value = "Liliana Fuchs"
$ws.Range("D21").Value = ""
$ws.Range("E21").Value = ""
$ws.Range("F21").Value = "https://www.directoalpaladar.com/recetas-vegetarianas/no-sabia-que-mis-ensaladas-legumbres-tenian-fallo-todo-ha-cambiado-que-no-me-salto-este-truco"

# Re-create hyperlinks for column F
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.directoalpaladar.com/recetas-vegetarianas/lentejas-tambien-para-verano-esta-ensalada-facil-saludable-rapida-muy-rica")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.directoalpaladar.com/recetas-vegetarianas/espaguetis-salsa-berenjena-tomate-queso-feta-receta-casera-deliciosa-que-querras-hacer-otra-vez")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.directoalpaladar.com/recetas-vegetarianas/hamburguesa-garbanzos-bote-legumbres-tienes-comida-bien-rica")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.directoalpaladar.com/recetas-vegetarianas/arroz-al-coco-truco-para-darle-sabor-al-cereal-todos-dias")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.directoalpaladar.com/recetas-vegetarianas/palomitas-maiz-al-curry-picoteo-diferente-adictivo")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.directoalpaladar.com/recetas-vegetarianas/tallarines-pepino-yogur-receta-que-simplemente-grita-verano")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.directoalpaladar.com/recetas-vegetarianas/recetas-frescas-para-disfrutar-verano-menu-semanal-vegetariano-7-agosto")
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.directoalpaladar.com/recetas-vegetarianas/paradoja-queso-que-vegetarianos-no-quieren-comer-carne-les-dan-igual-lacteos")
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.directoalpaladar.com/recetas-vegetarianas/pate-rosa-remolacha-tofu-receta-para-decorar-mesa")
$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.directoalpaladar.com/recetas-vegetarianas/este-desayuno-huevos-tipico-oriente-medio-cena-veraniega-perfecta")
$ws.Hyperlinks.Add($ws.Range("F12"), "https://www.directoalpaladar.com/recetas-vegetarianas/esta-comida-callejera-italiana-sabrosa-facil-mi-receta-favorita-calabacin")
$ws.Hyperlinks.Add($ws.Range("F13"), "https://www.directoalpaladar.com/recetas-vegetarianas/muffins-salados-boniato-receta-para-chuparse-dedos")
$ws.Hyperlinks.Add($ws.Range("F14"), "https://www.directoalpaladar.com/recetas-vegetarianas/refrescante-saciante-encender-fuego-esta-receta-calabacin-cena-saludable-verano-perfecta")
$ws.Hyperlinks.Add($ws.Range("F15"), "https://www.directoalpaladar.com/recetas-vegetarianas/asi-nuevos-postres-veganos-dia-que-no-yogures-tres-bases-vegetales-diferentes")
$ws.Hyperlinks.Add($ws.Range("F16"), "https://www.directoalpaladar.com/recetas-vegetarianas/despedimos-septimo-mes-ano-menu-semanal-vegetariano-31-julio")
$ws.Hyperlinks.Add($ws.Range("F17"), "https://www.directoalpaladar.com/recetas-vegetarianas/palitos-queso-receta-facil-rapida-aperitivo-verano")
$ws.Hyperlinks.Add($ws.Range("F18"), "https://www.directoalpaladar.com/recetas-vegetarianas/galletas-naranja-harina-garbanzos-avena-receta-facil-para-acompanar-merienda")
$ws.Hyperlinks.Add($ws.Range("F19"), "https://www.directoalpaladar.com/recetas-vegetarianas/receta-ensalada-coreana-patata-gamja-cremosa-refrescante-que-resulta-adictiva")
$ws.Hyperlinks.Add($ws.Range("F20"), "https://www.directoalpaladar.com/recetas-vegetarianas/recetas-ligerisimas-para-disfrutar-verano-menu-semanal-vegetariano-24-julio")
$ws.Hyperlinks.Add($ws.Range("F21"), "https://www.directoalpaladar.com/recetas-vegetarianas/no-sabia-que-mis-ensaladas-legumbres-tenian-fallo-todo-ha-cambiado-que-no-me-salto-este-truco")

# Restore hyperlink cell style on column F
$ws.Range("F2:F21").Style = "Hyperlink"